$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03171436517933
$ws.Range("D2").Value = 1.055960129672362
$ws.Range("E2").Value = 1.031250237884396
$ws.Range("F2").Value = 1.059412329378979
$ws.Range("I2").Value = 1.044126946388541
$ws.Range("J2").Value = 1.036848569362812
$ws.Range("K2").Value = 1.058699152972836
$ws.Range("L2").Value = 1.034058899110663
$ws.Range("M2").Value = 1.062141903793278
$ws.Range("N2").Value = 1.0383210132595

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032626417823843
$ws.Range("D3").Value = 1.05671118962624
$ws.Range("E3").Value = 1.032024219988418
$ws.Range("F3").Value = 1.060338388259914
$ws.Range("I3").Value = 1.044393710848061
$ws.Range("J3").Value = 1.037402803485467
$ws.Range("K3").Value = 1.059263916739493
$ws.Range("L3").Value = 1.034641498399953
$ws.Range("M3").Value = 1.062881904074405
$ws.Range("N3").Value = 1.038876034458181

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033217045858615
$ws.Range("D4").Value = 1.05719727443511
$ws.Range("E4").Value = 1.032525818448516
$ws.Range("F4").Value = 1.060938235728272
$ws.Range("I4").Value = 1.04456511889427
$ws.Range("J4").Value = 1.037761274941911
$ws.Range("K4").Value = 1.059628789410836
$ws.Range("L4").Value = 1.035018600861784
$ws.Range("M4").Value = 1.063360728071007
$ws.Range("N4").Value = 1.03923501498521

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033465457082592
$ws.Range("D5").Value = 1.057401646623391
$ws.Range("E5").Value = 1.032736875826509
$ws.Range("F5").Value = 1.061190559814199
$ws.Range("I5").Value = 1.044636889337181
$ws.Range("J5").Value = 1.037911938242098
$ws.Range("K5").Value = 1.059782044890492
$ws.Range("L5").Value = 1.035177162820825
$ws.Range("M5").Value = 1.063562022984829
$ws.Range("N5").Value = 1.039385892244565

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033507172909616
$ws.Range("D6").Value = 1.057435962866928
$ws.Range("E6").Value = 1.032772324118876
$ws.Range("F6").Value = 1.061232934785349
$ws.Range("I6").Value = 1.044648922914034
$ws.Range("J6").Value = 1.037937233010664
$ws.Range("K6").Value = 1.05980776909203
$ws.Range("L6").Value = 1.035203787687318
$ws.Range("M6").Value = 1.063595821089726
$ws.Range("N6").Value = 1.039411222934604

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033220364705956
$ws.Range("D7").Value = 1.057200005181225
$ws.Range("E7").Value = 1.032528637880872
$ws.Range("F7").Value = 1.060941606714603
$ws.Range("I7").Value = 1.044566079032833
$ws.Range("J7").Value = 1.037763288263228
$ws.Range("K7").Value = 1.059630837757833
$ws.Range("L7").Value = 1.035020719465428
$ws.Range("M7").Value = 1.063363417795762
$ws.Range("N7").Value = 1.039237031165675

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032022499884766
$ws.Range("D8").Value = 1.056213932531792
$ws.Range("E8").Value = 1.031511646345997
$ws.Range("F8").Value = 1.059725164894339
$ws.Range("I8").Value = 1.044217350066567
$ws.Range("J8").Value = 1.037035906904836
$ws.Range("K8").Value = 1.058890134269323
$ws.Range("L8").Value = 1.034255765414866
$ws.Range("M8").Value = 1.062391990722374
$ws.Range("N8").Value = 1.038508616842322

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029915348521731
$ws.Range("D9").Value = 1.054477170137107
$ws.Range("E9").Value = 1.029725615884296
$ws.Range("F9").Value = 1.05758649641129
$ws.Range("I9").Value = 1.043593635781981
$ws.Range("J9").Value = 1.03575302322711
$ws.Range("K9").Value = 1.057580628137342
$ws.Range("L9").Value = 1.032908798283875
$ws.Range("M9").Value = 1.06068022770977
$ws.Range("N9").Value = 1.037223911322625

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028513091030518
$ws.Range("D10").Value = 1.053319971892096
$ws.Range("E10").Value = 1.028539072687859
$ws.Range("F10").Value = 1.056164073102792
$ws.Range("I10").Value = 1.043171674817877
$ws.Range("J10").Value = 1.034897055047055
$ws.Range("K10").Value = 1.056704808559048
$ws.Range("L10").Value = 1.032011541798521
$ws.Range("M10").Value = 1.059539141687827
$ws.Range("N10").Value = 1.036366727569573

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027906506352577
$ws.Range("D11").Value = 1.052819063566435
$ws.Range("E11").Value = 1.028026286752656
$ws.Range("F11").Value = 1.055548962708815
$ws.Range("I11").Value = 1.042987510373645
$ws.Range("J11").Value = 1.034526254286263
$ws.Range("K11").Value = 1.056324915805936
$ws.Range("L11").Value = 1.031623204096148
$ws.Range("M11").Value = 1.05904507641158
$ws.Range("N11").Value = 1.035995400229169

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027681285245589
$ws.Range("D12").Value = 1.052633030354927
$ws.Range("E12").Value = 1.027835966126188
$ws.Range("F12").Value = 1.055320606232607
$ws.Range("I12").Value = 1.042918885906125
$ws.Range("J12").Value = 1.03438849913038
$ws.Range("K12").Value = 1.056183709143565
$ws.Range("L12").Value = 1.031478986302048
$ws.Range("M12").Value = 1.058861564740585
$ws.Range("N12").Value = 1.035857449445165

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027729591783195
$ws.Range("D13").Value = 1.052672933889901
$ws.Range("E13").Value = 1.027876783695548
$ws.Range("F13").Value = 1.055369583886198
$ws.Range("I13").Value = 1.042933615923642
$ws.Range("J13").Value = 1.034418049115834
$ws.Range("K13").Value = 1.056214002849303
$ws.Range("L13").Value = 1.031509920219368
$ws.Range("M13").Value = 1.058900928317031
$ws.Range("N13").Value = 1.035887041394987

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027887887634271
$ws.Range("D14").Value = 1.052803685455135
$ws.Range("E14").Value = 1.028010551703863
$ws.Range("F14").Value = 1.055530084175351
$ws.Range("I14").Value = 1.042981842288888
$ws.Range("J14").Value = 1.03451486787089
$ws.Range("K14").Value = 1.056313245603154
$ws.Range("L14").Value = 1.031611282428688
$ws.Range("M14").Value = 1.059029907146275
$ws.Range("N14").Value = 1.03598399764378

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027985431080517
$ws.Range("D15").Value = 1.052884249358143
$ws.Range("E15").Value = 1.028092990626465
$ws.Range("F15").Value = 1.055628990031865
$ws.Range("I15").Value = 1.043011527330115
$ws.Range("J15").Value = 1.034574518049115
$ws.Range("K15").Value = 1.056374379443209
$ws.Range("L15").Value = 1.031673738788262
$ws.Range("M15").Value = 1.059109376116045
$ws.Range("N15").Value = 1.036043732532101

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028553360827722
$ws.Range("D16").Value = 1.053353219155044
$ws.Range("E16").Value = 1.02857312568661
$ws.Range("F16").Value = 1.056204913109455
$ws.Range("I16").Value = 1.04318386665612
$ws.Range("J16").Value = 1.034921660561809
$ws.Range("K16").Value = 1.056730007043152
$ws.Range("L16").Value = 1.032037318379414
$ws.Range("M16").Value = 1.059571931967734
$ws.Range("N16").Value = 1.036391368026981

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028909770046847
$ws.Range("D17").Value = 1.053647437048038
$ws.Range("E17").Value = 1.028874569129164
$ws.Range("F17").Value = 1.056566391883196
$ws.Range("I17").Value = 1.043291582119995
$ws.Range("J17").Value = 1.035139371468972
$ws.Range("K17").Value = 1.056952907724018
$ws.Range("L17").Value = 1.032265431246949
$ws.Range("M17").Value = 1.059862090643531
$ws.Range("N17").Value = 1.036609388108603

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029117715538302
$ws.Range("D18").Value = 1.053819065299121
$ws.Range("E18").Value = 1.029050491908005
$ws.Range("F18").Value = 1.056777314224958
$ws.Range("I18").Value = 1.043354270533607
$ws.Range("J18").Value = 1.035266342971846
$ws.Range("K18").Value = 1.057082858443137
$ws.Range("L18").Value = 1.032398502947913
$ws.Range("M18").Value = 1.060031338339536
$ws.Range("N18").Value = 1.036736539925576

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.0291886293887
$ws.Range("D19").Value = 1.053877588754707
$ws.Range("E19").Value = 1.029110493256249
$ws.Range("F19").Value = 1.056849246459163
$ws.Range("I19").Value = 1.0433756218639
$ws.Range("J19").Value = 1.035309634301216
$ws.Range("K19").Value = 1.057127157478609
$ws.Range("L19").Value = 1.032443879866153
$ws.Range("M19").Value = 1.060089047908197
$ws.Range("N19").Value = 1.036779892733599

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028871524701745
$ws.Range("D20").Value = 1.053615868593943
$ws.Range("E20").Value = 1.028842217169157
$ws.Range("F20").Value = 1.056527600564687
$ws.Range("I20").Value = 1.043280039766439
$ws.Range("J20").Value = 1.035116014754426
$ws.Range("K20").Value = 1.056928999147912
$ws.Range("L20").Value = 1.032240955095638
$ws.Range("M20").Value = 1.059830959038903
$ws.Range("N20").Value = 1.036585998224844

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027841270919042
$ws.Range("D21").Value = 1.052765181636209
$ws.Range("E21").Value = 1.027971156182637
$ws.Range("F21").Value = 1.055482817435733
$ws.Range("I21").Value = 1.042967646827835
$ws.Range("J21").Value = 1.034486357794841
$ws.Range("K21").Value = 1.056284023770195
$ws.Range("L21").Value = 1.03158143299989
$ws.Range("M21").Value = 1.058991925915085
$ws.Range("N21").Value = 1.035955447080154

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027194039076573
$ws.Range("D22").Value = 1.052230474709083
$ws.Range("E22").Value = 1.027424359330692
$ws.Range("F22").Value = 1.054826632676863
$ws.Range("I22").Value = 1.042769974379427
$ws.Range("J22").Value = 1.034090333454428
$ws.Range("K22").Value = 1.055877938153913
$ws.Range("L22").Value = 1.031166928483916
$ws.Range("M22").Value = 1.058464428829654
$ws.Range("N22").Value = 1.035558860339752

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027537098388327
$ws.Range("D23").Value = 1.052513917902294
$ws.Range("E23").Value = 1.027714143398009
$ws.Range("F23").Value = 1.05517442066362
$ws.Range("I23").Value = 1.042874883372558
$ws.Range("J23").Value = 1.03430028584346
$ws.Range("K23").Value = 1.056093264901523
$ws.Range("L23").Value = 1.031386649355285
$ws.Range("M23").Value = 1.058744061140359
$ws.Range("N23").Value = 1.035769110885259

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02888880594773
$ws.Range("D24").Value = 1.053630132970427
$ws.Range("E24").Value = 1.02885683533109
$ws.Range("F24").Value = 1.056545128451249
$ws.Range("I24").Value = 1.043285255692404
$ws.Range("J24").Value = 1.035126568696593
$ws.Range("K24").Value = 1.056939802600563
$ws.Range("L24").Value = 1.032252014761424
$ws.Range("M24").Value = 1.059845026061689
$ws.Range("N24").Value = 1.036596567154819

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030459659870839
$ws.Range("D25").Value = 1.054926058006302
$ws.Range("E25").Value = 1.03018662333195
$ws.Range("F25").Value = 1.058138808593565
$ws.Range("I25").Value = 1.043755968068203
$ws.Range("J25").Value = 1.036084809170175
$ws.Range("K25").Value = 1.057919668247408
$ws.Range("L25").Value = 1.033256898995281
$ws.Range("M25").Value = 1.061122749507856
$ws.Range("N25").Value = 1.037556168439785

